$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: value removed -> becomes a present-but-blank (empty text) cell
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"

# C4: value changed to 0
$ws.Range("C4").Value = 0

# C5: value changed
$ws.Range("C5").Value = 619.3421170982775

# Row 7 ("Other") relabeled to "Biogas" with a new D7 value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 656.6793065786602

# New row 8 ("Other"): copy A7's formatting (bold/border label style) onto A8,
# then set its text; B8/C8 stay blank (empty text) but present; D8 gets the new value.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = "'"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = 10100.5669154866
